$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the table (Task, Mean, Std, Seed_0, Seed_1) - Seed_2 column (F) is left blank
$rows = @(
    @("Task 1.1",  4449.200320822603, 388.4662306475441, 4060.734090175059, 4837.666551470147),
    @("Task 1.2",  4449.200320822603, 388.4662306475441, 4060.734090175059, 4837.666551470147),
    @("Task 1.3",  5079.337457256514, 571.3560627366455, 5650.693519993159, 4507.981394519868),
    @("Task 2.1",  18290.62894347315, 6147.148928657862, 12143.48001481529, 24437.77787213101),
    @("Task 2.2",  18290.62894347315, 6147.148928657862, 12143.48001481529, 24437.77787213101),
    @("Task 3.1.1", 4192.381535149367, 3017.902992305311, 7210.284527454678, 1174.478542844056)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = ""
    $r = $r + 1
}
